$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05354133333333333
$ws.Range("H2").Value = 0.160624
$ws.Range("I2").Value = 0.00209946492164722
$ws.Range("J2").Value = 0.00209946492164722
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.341267
$ws.Range("N2").Value = 7.023801
$ws.Range("Q2").Value = 0.1253545568693333
$ws.Range("R2").Value = 1.128191011824
$ws.Range("S2").Value = 0.00209946492164722
$ws.Range("T2").Value = 0.00209946492164722

# Row 3
$ws.Range("I3").Value = 0.05460670042535784
$ws.Range("J3").Value = 0.05460670042535784
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.341267
$ws.Range("N3").Value = 7.023801
$ws.Range("Q3").Value = 3.260449204622333
$ws.Range("R3").Value = 29.344042841601
$ws.Range("S3").Value = 0.05460670042535784
$ws.Range("T3").Value = 0.05460670042535784

# Row 4
$ws.Range("G4").Value = 24.05622933333333
$ws.Range("H4").Value = 72.168688
$ws.Range("I4").Value = 0.943293834652995
$ws.Range("J4").Value = 0.943293834652995
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.341267
$ws.Range("N4").Value = 7.023801
$ws.Range("Q4").Value = 56.32205588256533
$ws.Range("R4").Value = 506.898502943088
$ws.Range("S4").Value = 0.943293834652995
$ws.Range("T4").Value = 0.943293834652995
